# Apply updated graph data values and cell selection change

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 values (A6:E6)
$ws.Range("A6").Value = 2764.5780420000001
$ws.Range("B6").Value = 253.9402216
$ws.Range("C6").Value = 193.24799999999999
$ws.Range("D6").Value = 194.08235110000001
$ws.Range("E6").Value = 192.41388810000001

# Update a few individual B-column values elsewhere in the table
$ws.Range("B16").Value = 253.9402216
$ws.Range("B18").Value = 253.9402216
$ws.Range("B25").Value = 253.9402216
$ws.Range("B27").Value = 253.9402216

# Update the active selection to G8
$ws.Range("G8").Select()
